$d = $word.ActiveDocument

# Character 11 (vertical-tab) is how Word's Range.Text API represents a
# manual line break (<w:br/>) when you assign text containing one.
$brk = [char]11

# Helper: insert a brand-new paragraph immediately before the paragraph
# currently sitting at index $i (pushing it and everything after it down
# by one), optionally giving the new paragraph some text. Returns the
# index the (now shifted) original paragraph sits at, i.e. $i + 1, so the
# caller can keep chaining inserts in document order.
function Insert-ParaBeforeIndex($doc, $i, $text) {
    $anchor = $doc.Paragraphs.Item($i)
    $anchor.Range.InsertParagraphBefore()
    if ($text -ne "") {
        $newp = $doc.Paragraphs.Item($i)
        $newp.Range.Text = $text
    }
    return $i + 1
}

# --- 1. "Необходимо заполнить раздел..." -> prefix with "1. " ---
$d.Content.Find.Execute(
    "Необходимо заполнить раздел",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1. Необходимо заполнить раздел", 2) | Out-Null

# That paragraph is (still) #2. The next paragraph down is #4 (after the
# blank spacer #3). Insert the brand-new "2. ..." paragraph + its spacer
# right before the existing "Согласно п.5.9. ..." paragraph (#4).
$idx = 4
$idx = Insert-ParaBeforeIndex $d $idx "2. Обратите внимание на п.5.3. Методических рекомендаций: необходимо, чтобы на подтверждающих документах стояла подпись ответственного лица - руководителя практики или руководителя организации."
$idx = Insert-ParaBeforeIndex $d $idx ""

# --- 3. "Согласно п.5.9. ..." -> prefix with "3. " ---
# ($idx now points at the "Согласно п.5.9. ..." paragraph)
$d.Content.Find.Execute(
    "Согласно п.5.9. Методических рекомендаций",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "3. Согласно п.5.9. Методических рекомендаций", 2) | Out-Null

# Move past the now-prefixed paragraph #3 AND the blank spacer that
# already followed "Согласно п.5.9. ..." in the original document (that
# existing spacer stays put, sitting between #3 and the new #4). The next
# insertion point is therefore two slots down, right before the closing
# "Ждем доработанный..." paragraph.
$idx = $idx + 2

# --- 4..9: brand new paragraphs (each with its own blank spacer before
#     the closing paragraph) ---
$text4a = "4. Необходимо ознакомиться с разделом 5 Методических рекомендаций по подготовке Мониторингового отчета. К каждому показателю, по которому произошли изменения в отчетном периоде, необходимо приложить подтверждающий документ.  Такой документ должен быть оформлен соответствующим образом (см.п.5.4, 5.5, 5.8. Метод.рекомендаций) и скреплен подписью ответственного лица и печатью организации. Если изменений в значении показателя не было в отчетном периоде, нужно поставить значение 0 и написать пояснительный комментарий."
$text4b = "Обратите внимание, что по основным показателям (пок. 1, 2, 3 и т.д.) обязательно предоставлять подтверждающие документы; по подпоказателям (1.1, 2.1, 3.1.5 и т.д.) отдельные подтверждающие документы не обязательны, если сведения по ним могут быть отражены в подтверждающем документе к основному показателю."
$text4 = $text4a + $brk + $brk + $text4b
$idx = Insert-ParaBeforeIndex $d $idx $text4
$idx = Insert-ParaBeforeIndex $d $idx ""

$idx = Insert-ParaBeforeIndex $d $idx "5. В подтверждающем документе содержится информация о работе, проделанной вне дат отчетного периода. Просим исключить эти данные из подтверждающих документов и скорректировать значение показателя. "
$idx = Insert-ParaBeforeIndex $d $idx ""

$idx = Insert-ParaBeforeIndex $d $idx "6. Среди подтверждающих документов нет единого сводного документа, дающего представление о всей проделанной работе/оказанных услугах целевым группам. Необходимо сформировать сводную таблицу оказанных услуг согласно методическим рекомендациям. "
$idx = Insert-ParaBeforeIndex $d $idx ""

$idx = Insert-ParaBeforeIndex $d $idx "7. В соответствии с п.5.5. Методических рекомендаций необходимо корректно оформить подтверждающие документы. В “шапке” документа нужно указать название организации, отчетный период, формулировку показателя, значение показателя."
$idx = Insert-ParaBeforeIndex $d $idx ""

$idx = Insert-ParaBeforeIndex $d $idx "8. Обратите внимание, что значения показателя в мониторинговой форме и в подтверждающем документе не совпадают. Такого расхождения быть не должно. Необходимо привести в соответствие данные.  "
$idx = Insert-ParaBeforeIndex $d $idx ""

$idx = Insert-ParaBeforeIndex $d $idx "9. Мы не можем принять значение показателя к зачету, т.к. регистрация - некорректный в данном случае метод сбора данных. В качестве метода сбора данных, которые фиксируют изменения в жизни благополучателей, могут использоваться анкетирование, тестирование, диагностика, наблюдение и т.д. При составлении заявки вы указали метод ....."
